# Refresh the crypto price/volume snapshot (cryptos list update, GitHub Actions run).
# Price (D) and Volume(1h) (E) columns hold plain text in the source sheet; some
# "Price" strings are pure decimals (e.g. "5.58") which Excel would otherwise
# auto-coerce into numbers on assignment, so those are written with a leading
# apostrophe to force a text entry (matches the original General/text formatting,
# i.e. how this would be done interactively in Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.501.39"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.640.47"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D5").Value = "'603.70"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'146.28"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.109"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'5.58"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  +4.66%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'27.53"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "3.116.65"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "63.298.54"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "2.652.92"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "'11.46"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'4.57"
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("D20").Value = "'343.27"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +3.22%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "'1.70"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "'9.09"
$ws.Range("E26").Value = "  +8.44%  "
$ws.Range("D27").Value = "'577.96"
$ws.Range("E27").Value = "  +6.53%  "
$ws.Range("D28").Value = "'1.56"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").Value = "0.0₃0823"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").Value = "  +7.29%  "
$ws.Range("D36").Value = "'166.96"
$ws.Range("E36").Value = "  -4.66%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'1.94"
$ws.Range("E39").Value = "  +7.99%  "
$ws.Range("D40").Value = "'19.13"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'168.86"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").Value = "'22.14"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "'0.0570"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'1.90"
$ws.Range("E49").Value = "  +11.87%  "
$ws.Range("D50").Value = "'18.76"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").Value = "'11.27"
$ws.Range("E51").Value = "  -0.67%  "
